# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5485
$wsExpo.Range("F3").Value = 108
$wsExpo.Range("F4").Value = 348
$wsExpo.Range("F7").Value = 52
$wsExpo.Range("F11").Value = 420
$wsExpo.Range("F12").Value = 2999
$wsExpo.Range("F13").Value = 163
$wsExpo.Range("F14").Value = 1602

# --- Sheet: 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5485
$wsAll.Range("F3").Value = 108
$wsAll.Range("F4").Value = 348
$wsAll.Range("F8").Value = 52
$wsAll.Range("F12").Value = 420
$wsAll.Range("F13").Value = 2999
$wsAll.Range("F14").Value = 163
$wsAll.Range("F15").Value = 1602
